$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the existing data row (row 2), shifting the
# current rows 2-3 (K M Swaminadh / ADMIN) down to rows 4-5, and bring their
# new neighbours' formatting along for the ride.
$ws.Rows("2:3").Insert()

# Fill in the new row 3 (Sakshi Jain) and row 2 (Ammiraju Rajasekhar) -
# userid for both of them is the placeholder "AExxxxx".
$ws.Range("A3").Value = "Sakshi Jain"
$ws.Range("C3").Value = "sakshi_jain@goodyear.com"
$ws.Range("A2").Value = "Ammiraju Rajasekhar"
$ws.Range("C2").Value = "ammiraju_rajasekhar@goodyear.com"
$ws.Range("B2").Value = "AExxxxx"
$ws.Range("B3").Value = "AExxxxx"

# The two original mailto hyperlinks (now anchored on C4/C5 after the row
# insert) kept pointing at their old C2/C3 addresses, so re-create them on
# the correct, shifted-down cells.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:swaminadh_kone@goodyear.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:swamynadh.engg@gmail.com")

# Hyperlinks.Add() re-stamps the anchor cell with a brand new (duplicate)
# "Hyperlink" style; reapply the worksheet's existing named Hyperlink style
# so the cells keep using the original style record instead of a clone.
$ws.Range("C4").Style = "Hyperlink"
$ws.Range("C5").Style = "Hyperlink"

$ws.Range("A2").Select()
